# This workbook has two near-identical "原始凭证" (original voucher) sheets,
# each driven by a handful of raw input cells (total amount, rate table,
# on-time timestamp, and session-length minutes) with everything else
# (off-time, usage fee, platform fee, card benefit, consumption, balance)
# computed by formulas that key off those inputs via named ranges.
#
# The edit only changes the raw inputs; the formulas (and their text) are
# left untouched so the workbook keeps computing the dependent figures.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1 (tab 1): 浙江杭州滨江中南乐游城店_原始凭证 ----
$ws1 = $wb.Worksheets.Item(1)

# 总金额 (total amount)
$ws1.Range("F1").Value2 = 74.63

# 区域费率_高级_假期 row of the C6:G9 lookup table
$ws1.Range("C9").Value2 = 5
$ws1.Range("D9").Value2 = 10
$ws1.Range("E9").Value2 = 13
$ws1.Range("F9").Value2 = 16
$ws1.Range("G9").Value2 = 17

# 上机时间 (session start time)
$ws1.Range("B14").Value2 = 45471.519618055558

# 上机时长_分 (session length - minutes component)
$ws1.Range("F15").Value2 = 2

# ---- Sheet 2 (tab 2): 浙江杭州西湖三墩地铁站店_原始凭证 ----
$ws2 = $wb.Worksheets.Item(2)

# 总金额 (total amount)
$ws2.Range("F1").Value2 = 74.63

# 上机时间 (session start time)
$ws2.Range("B14").Value2 = 45471.519618055558

# 上机时长_分 (session length - minutes component)
$ws2.Range("F15").Value2 = 2
